$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "members" column (C) with new community membership lists
$ws.Range("C2").Value = "['Fc5', 'Fc3', 'Fc1', 'C5', 'Cp5', 'Fp1', 'Fpz', 'Af7', 'Af3', 'F7', 'F5', 'F3', 'Ft7', 'T7', 'Tp7']"
$ws.Range("C3").Value = "['Fcz', 'Fc2', 'C1', 'Cz', 'C2', 'Afz', 'F1', 'Fz', 'F2']"
$ws.Range("C4").Value = "['Fc4', 'Fc6', 'C4', 'C6', 'Fp2', 'Af4', 'Af8', 'F4', 'F6', 'F8', 'Ft8', 'T8']"
$ws.Range("C5").Value = "['C3', 'Cp3', 'Cp1', 'T9', 'T10', 'P7', 'P5', 'P3', 'P1', 'Po7', 'Po3', 'O1', 'Iz']"
$ws.Range("C6").Value = "['Cpz', 'Cp2', 'Cp4', 'Cp6', 'Tp8', 'Pz', 'P2', 'P4', 'P6', 'P8', 'Poz', 'Po4', 'Po8', 'Oz', 'O2']"

# Update "community" column (B) values
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 2
